$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "37.213.31"
$ws.Range("E2").Value = "  +1.69%  "

# Row 3
$ws.Range("D3").Value = "2.058.51"
$ws.Range("E3").Value = "  +1.18%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.16"

# Row 7
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("E8").Value = "  +3.47%  "

# Row 9
$ws.Range("E9").Value = "  +3.11%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.89"
$ws.Range("E10").Value = "  +0.63%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0758"

# Row 12
$ws.Range("E12").Value = "  +1.39%  "

# Row 13
$ws.Range("D13").Value = "2.360.22"
$ws.Range("E13").Value = "  +0.97%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.53"
$ws.Range("E14").Value = "  +1.85%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.79"
$ws.Range("E15").Value = "  +3.79%  "

# Row 16
$ws.Range("E16").Value = "  +2.27%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.15"
$ws.Range("E17").Value = "  -0.20%  "

# Row 18
$ws.Range("D18").Value = "2.055.99"
$ws.Range("E18").Value = "  +1.05%  "

# Row 19
$ws.Range("D19").Value = "37.148.58"
$ws.Range("E19").Value = "  +1.02%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.39"
$ws.Range("E20").Value = "  +9.38%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.19"
$ws.Range("E21").Value = "  +2.21%  "

# Row 22
$ws.Range("E22").Value = "  +1.54%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "225.91"
$ws.Range("E23").Value = "  +2.35%  "

# Row 24
$ws.Range("E24").Value = "  -0.09%  "

# Row 25
$ws.Range("E25").Value = "  +0.46%  "

# Row 26
$ws.Range("E26").Value = "  +0.83%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.50"
$ws.Range("E27").Value = "  +1.68%  "

# Row 28
$ws.Range("E28").Value = "  +7.69%  "

# Row 29
$ws.Range("E29").Value = "  +0.53%  "

# Row 30
$ws.Range("E30").Value = "  +0.68%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.04"
$ws.Range("E31").Value = "  +0.77%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.117"
$ws.Range("E32").Value = "  +0.35%  "

# Row 33
$ws.Range("E33").Value = "  +2.16%  "

# Row 34
$ws.Range("E34").Value = "  +1.72%  "

# Row 35
$ws.Range("E35").Value = "  +7.73%  "

# Row 36
$ws.Range("E36").Value = "  +1.10%  "

# Row 37
$ws.Range("E37").Value = "  +0.00%  "

# Row 38
$ws.Range("E38").Value = "  -0.55%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.25"
$ws.Range("E39").Value = "  +1.24%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.67"
$ws.Range("E40").Value = "  -1.74%  "

# Row 41
$ws.Range("E41").Value = "  +0.07%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.39"
$ws.Range("E42").Value = "  -1.71%  "

# Row 43
$ws.Range("D43").Value = "1.468.75"
$ws.Range("E43").Value = "  -0.57%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.31"
$ws.Range("E44").Value = "  +2.75%  "

# Row 45
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0931"
$ws.Range("E45").Value = "  -1.56%  "

# Row 46
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.17"
$ws.Range("E46").Value = "  +5.33%  "

# Row 47
$ws.Range("E47").Value = "  +3.63%  "

# Row 48
$ws.Range("E48").Value = "  +1.62%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.06"
$ws.Range("E49").Value = "  -3.42%  "

# Row 50
$ws.Range("E50").Value = "  +2.78%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.95"
$ws.Range("E51").Value = "  +1.68%  "
